# Scheduled-runner price/profit refresh across the per-job (ALC/ARM/BSM/CRP/
# CUL/GSM/LTW/WVR) leve-profit sheets: updates currentAveragePrice* /
# LevePrice* / LeveProfit* columns (H-N) for the rows whose market prices
# moved since the last run.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 747909.3
$ws.Range("J17").Value = 747909.3
$ws.Range("L17").Value = 2243727.9
$ws.Range("N17").Value = -2244063.9
$ws.Range("H29").Value = 793.6
$ws.Range("I29").Value = 793.6
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 2380.8
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -2099.8
$ws.Range("N29").ClearContents()
$ws.Range("H38").Value = 935.7143
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 935.7143
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 2807.1429
$ws.Range("M38").ClearContents()
$ws.Range("N38").Value = -3551.1429
$ws.Range("H58").Value = 1000.25
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 1000.25
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 3000.75
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -3300.75
$ws.Range("H61").Value = 13371.25
$ws.Range("I61").Value = 16745
$ws.Range("K61").Value = 50235
$ws.Range("M61").Value = -50063
$ws.Range("H112").Value = 1346.6666
$ws.Range("I112").Value = 853.3333
$ws.Range("K112").Value = 2559.9999
$ws.Range("M112").Value = -1451.9999
$ws.Range("H137").Value = 2985.6667
$ws.Range("I137").Value = 1265.7646
$ws.Range("J137").Value = 4656.4287
$ws.Range("K137").Value = 3797.2938
$ws.Range("L137").Value = 13969.2861
$ws.Range("M137").Value = -1247.2938
$ws.Range("N137").Value = -19069.2861

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1562.65
$ws.Range("I61").Value = 1384.0322
$ws.Range("K61").Value = 1384.0322
$ws.Range("M61").Value = -1172.0322
$ws.Range("H74").Value = 3799.647
$ws.Range("I74").Value = 788.6667
$ws.Range("J74").Value = 26382
$ws.Range("K74").Value = 788.6667
$ws.Range("L74").Value = 26382
$ws.Range("M74").Value = 85.33330000000001
$ws.Range("N74").Value = -28130
$ws.Range("H77").Value = 3799.647
$ws.Range("I77").Value = 788.6667
$ws.Range("J77").Value = 26382
$ws.Range("K77").Value = 3943.3335
$ws.Range("L77").Value = 131910
$ws.Range("M77").Value = 424.6665000000003
$ws.Range("N77").Value = -140646
$ws.Range("H136").Value = 1562.65
$ws.Range("I136").Value = 1384.0322
$ws.Range("K136").Value = 4152.096600000001
$ws.Range("M136").Value = -1602.096600000001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 46517.684
$ws.Range("I20").Value = 1070.2858
$ws.Range("J20").Value = 126050.625
$ws.Range("K20").Value = 1070.2858
$ws.Range("L20").Value = 126050.625
$ws.Range("M20").Value = -823.2858000000001
$ws.Range("N20").Value = -126544.625
$ws.Range("H134").Value = 1099.8918
$ws.Range("I134").Value = 1003.06665
$ws.Range("K134").Value = 3009.19995
$ws.Range("M134").Value = -474.1999500000002

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 27282.057
$ws.Range("I31").Value = 2904.6365
$ws.Range("K31").Value = 2904.6365
$ws.Range("M31").Value = -2609.6365
$ws.Range("H34").Value = 27282.057
$ws.Range("I34").Value = 2904.6365
$ws.Range("K34").Value = 2904.6365
$ws.Range("M34").Value = -2702.6365
$ws.Range("H58").Value = 863.45764
$ws.Range("I58").Value = 767.6667
$ws.Range("J58").Value = 1281.4546
$ws.Range("K58").Value = 767.6667
$ws.Range("L58").Value = 1281.4546
$ws.Range("M58").Value = -564.6667
$ws.Range("N58").Value = -1687.4546
$ws.Range("H132").Value = 13892108
$ws.Range("I132").Value = 20411916
$ws.Range("K132").Value = 61235748
$ws.Range("M132").Value = -61233218
$ws.Range("H134").Value = 4002.487
$ws.Range("I134").Value = 4786.6333
$ws.Range("J134").Value = 1388.6666
$ws.Range("K134").Value = 14359.8999
$ws.Range("L134").Value = 4165.9998
$ws.Range("M134").Value = -11824.8999
$ws.Range("N134").Value = -9235.9998
$ws.Range("H136").Value = 863.45764
$ws.Range("I136").Value = 767.6667
$ws.Range("J136").Value = 1281.4546
$ws.Range("K136").Value = 2303.0001
$ws.Range("L136").Value = 3844.3638
$ws.Range("M136").Value = 246.9998999999998
$ws.Range("N136").Value = -8944.3638

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H121").Value = 19000
$ws.Range("J121").Value = 19000
$ws.Range("L121").Value = 19000
$ws.Range("N121").Value = -22494

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1628.909
$ws.Range("I68").Value = 1568.8182
$ws.Range("J68").Value = 1749.091
$ws.Range("K68").Value = 1568.8182
$ws.Range("L68").Value = 1749.091
$ws.Range("M68").Value = -819.8181999999999
$ws.Range("N68").Value = -3247.091
$ws.Range("H71").Value = 1628.909
$ws.Range("I71").Value = 1568.8182
$ws.Range("J71").Value = 1749.091
$ws.Range("K71").Value = 7844.090999999999
$ws.Range("L71").Value = 8745.455
$ws.Range("M71").Value = -4100.090999999999
$ws.Range("N71").Value = -16233.455
$ws.Range("H132").Value = 2979.4824
$ws.Range("I132").Value = 3334.776
$ws.Range("J132").Value = 2216.2593
$ws.Range("K132").Value = 10004.328
$ws.Range("L132").Value = 6648.777900000001
$ws.Range("M132").Value = -7474.328
$ws.Range("N132").Value = -11708.7779
$ws.Range("H136").Value = 2691.1638
$ws.Range("I136").Value = 966.65955
$ws.Range("J136").Value = 8480.571
$ws.Range("K136").Value = 2899.97865
$ws.Range("L136").Value = 25441.713
$ws.Range("M136").Value = -349.97865
$ws.Range("N136").Value = -30541.713
$ws.Range("H139").Value = 52000
$ws.Range("J139").Value = 52000
$ws.Range("L139").Value = 52000
$ws.Range("N139").Value = -62280
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()
$ws.Range("H141").Value = 133245
$ws.Range("J141").Value = 133245
$ws.Range("L141").Value = 133245
$ws.Range("N141").Value = -143605

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 18966.666
$ws.Range("J54").Value = 18966.666
$ws.Range("L54").Value = 18966.666
$ws.Range("N54").Value = -20006.666
$ws.Range("H132").Value = 2934.9678
$ws.Range("I132").Value = 3578.475
$ws.Range("J132").Value = 1764.9546
$ws.Range("K132").Value = 10735.425
$ws.Range("L132").Value = 5294.8638
$ws.Range("M132").Value = -8205.425
$ws.Range("N132").Value = -10354.8638
